$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-03 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-04 Friday", 2)

$d.Content.Find.Execute("274÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "585÷6=", 2)
$d.Content.Find.Execute("687÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "256÷8=", 2)
$d.Content.Find.Execute("777÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "335÷3=", 2)
$d.Content.Find.Execute("786÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "807÷7=", 2)
$d.Content.Find.Execute("569÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "735÷9=", 2)
$d.Content.Find.Execute("876÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "427÷4=", 2)
$d.Content.Find.Execute("435÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "953÷9=", 2)
$d.Content.Find.Execute("970÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "114÷6=", 2)
$d.Content.Find.Execute("742÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "215÷2=", 2)
$d.Content.Find.Execute("860÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "301÷7=", 2)
$d.Content.Find.Execute("359÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "210÷2=", 2)
$d.Content.Find.Execute("979÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "756÷2=", 2)
$d.Content.Find.Execute("950÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "803÷4=", 2)
$d.Content.Find.Execute("620÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "138÷6=", 2)
$d.Content.Find.Execute("585÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "752÷8=", 2)
$d.Content.Find.Execute("553÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "804÷8=", 2)
$d.Content.Find.Execute("349÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "256÷9=", 2)
$d.Content.Find.Execute("568÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "366÷9=", 2)
$d.Content.Find.Execute("346÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "682÷5=", 2)
$d.Content.Find.Execute("298÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "821÷2=", 2)
$d.Content.Find.Execute("808÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "703÷3=", 2)
$d.Content.Find.Execute("833÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "806÷6=", 2)
$d.Content.Find.Execute("943÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "223÷6=", 2)
$d.Content.Find.Execute("961÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "112÷8=", 2)
$d.Content.Find.Execute("773÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "909÷5=", 2)
